$d = $word.ActiveDocument
$d.Content.Find.Execute("neighbourhood", $true, $false, $false, $false, $false,
                         $true, 1, $false, "neighborhood", 2)
